$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: replace a paragraph's content via InsertXML, wrapping the supplied
# paragraph-body fragment ("$innerXml" - the <w:pPr>...</w:pPr> + runs) in a
# minimal OOXML package so Word's InsertXML applies it to the target Range.
# ---------------------------------------------------------------------------
function Set-ParaXml($range, [string]$innerXml) {
    $pkg = "<?xml version='1.0'?><pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body>" + $innerXml + "</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
    $range.InsertXML($pkg)
}

# ---------------------------------------------------------------------------
# 1) Paragraph 3: "Using bucki's availability..." -> split into multiple
#    runs with spellcheck/grammar proofErr markers around "bucki's" and
#    "memory".
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs(3).Range
Set-ParaXml $p3 @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Using </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>bucki&#8217;s</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> availability of 16 GB </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>memory</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>, computer difference of co-expression network.</w:t></w:r></w:p>
'@

# ---------------------------------------------------------------------------
# 2) Paragraph 5: "laziness" -> wrap in gramStart/gramEnd proofErr markers.
# ---------------------------------------------------------------------------
$p5 = $d.Paragraphs(5).Range
Set-ParaXml $p5 @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:t>laziness</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>
'@

# ---------------------------------------------------------------------------
# 3) Paragraph 8: "Using bucki's availability..." -> "Same as above".
# ---------------------------------------------------------------------------
$p8 = $d.Paragraphs(8).Range
Set-ParaXml $p8 @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Same as above</w:t></w:r></w:p>
'@

# ---------------------------------------------------------------------------
# 4) Paragraph 11: merge the two runs (previously split by a _GoBack
#    bookmark) into a single run with the complete sentence; the bookmark
#    is re-created later, at the end of the new content added below.
# ---------------------------------------------------------------------------
$p11 = $d.Paragraphs(11).Range
Set-ParaXml $p11 @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Need to streamline histogram with input parameters. Input false normalization flags could break the code.</w:t></w:r></w:p>
'@

# ---------------------------------------------------------------------------
# 5) Append the new "Tuesday October 7th" / "Wednesday October 8th" log
#    entries after paragraph 11.
# ---------------------------------------------------------------------------
$anchor = $d.Paragraphs(11).Range
$anchor.InsertParagraphAfter()

$newParasXml = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Tuesday October 7</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>th</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Goal-</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Same as above</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Expand R code to calculate spearman / other correlation metrics</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Roadblocks-</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Installing R packages locally on </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>bucki</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. Waiting for </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>oschelp</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> to install latest R version</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Wednesday October 8</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>th</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Goal</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Same as above</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>COMPLETED- Expand R code to calculate spearman</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Roadblocks-</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">OSC still working on </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Bucki</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> cluster.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@

$newAnchor = $d.Paragraphs(12).Range
Set-ParaXml $newAnchor $newParasXml

Write-Host "done"
